# Fix: Proxy and wirte excel
# Write the two report rows returned by the (now-fixed) proxy/API call into
# the worksheet, right below the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phone numbers are long, purely-numeric strings (e.g. "855764049358").
# Left as plain numbers Excel would happily mangle them (round-trip through
# floating point / scientific notation), so force column A to Text before
# writing the values - this is the actual bug the commit is fixing.
$ws.Range("A2:A3").NumberFormat = "@"

$data = @(
    @("855764049358", "success", "failed", "Message not found or Archived for another partner", "2024-07-25 22:17:07"),
    @("855764044995", "success", "failed", "Message not found or Archived for another partner", "2024-07-25 22:18:53")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $row = $row + 1
}

# The "Message" column needed to be widened so the longer text fits.
$ws.Range("D1").ColumnWidth = 50.14

# Mirror the author's final selection state over the newly written rows.
[void]$ws.Range("A2:E4").Select()
